$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 44313
$ws.Range("B3").Value = 71431
$ws.Range("C3").Value = 60000
$ws.Range("D3").Value = 3225
$ws.Range("E3").Value = 2019
$ws.Range("F3").Value = 1413
$ws.Range("G3").Value = 18666
$ws.Range("H3").Value = 1369
$ws.Range("I3").Value = 794
$ws.Range("J3").Value = 196
